# CryCompanywiseStockReport_1.xlsx edit
#
# The source report lists, for several SKUs, multiple consecutive rows that
# belong to the same item (same item name in column C, same rate in column D
# for most groups). For each such group of rows the values in columns
# B (barcode), D (MRP), E (rate), F (qty) and G (value) were cyclically
# rotated down by one row (the last row of the group wraps around to the
# first row), while columns A (serial no), C (item name) and H:M stay the
# same.
#
# This script re-creates that row rotation by writing the resulting
# (post-rotation) values directly into the affected cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-RowGroup($rows) {
    # $rows is an ordered list of row numbers that form one rotation group.
    # Columns that participate in the rotation:
    $cols = @(2, 4, 5, 6, 7)   # B, D, E, F, G

    # Snapshot current ("before") values for every cell that will move.
    $snapshot = @{}
    foreach ($r in $rows) {
        foreach ($c in $cols) {
            $snapshot["$r`_$c"] = $ws.Cells.Item($r, $c).Value2
        }
    }

    $n = $rows.Count
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $rows[$i]
        $srcRow  = $rows[($i - 1 + $n) % $n]
        foreach ($c in $cols) {
            $ws.Cells.Item($destRow, $c).Value = $snapshot["$srcRow`_$c"]
        }
    }
}

Rotate-RowGroup @(146, 147, 148)
Rotate-RowGroup @(277, 278)
Rotate-RowGroup @(292, 293)
Rotate-RowGroup @(294, 295, 296)
Rotate-RowGroup @(299, 300)
Rotate-RowGroup @(315, 316)
Rotate-RowGroup @(465, 466)
Rotate-RowGroup @(472, 473)
Rotate-RowGroup @(476, 477)
Rotate-RowGroup @(479, 480)
Rotate-RowGroup @(490, 491)
Rotate-RowGroup @(705, 706)
Rotate-RowGroup @(732, 733)
